$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column C: reciprocal of B (tau -> 1/tau) ---
$ws.Range("C2").Formula = "=1/B2"
$ws.Range("C3:C9").Formula = "=1/B3"

# --- New column J: (H-I)/(H+I) ---
$ws.Range("J2").Formula = "=(H2-I2)/(H2+I2)"
$ws.Range("J3:J7").Formula = "=(H3-I3)/(H3+I3)"

# --- New column K: reciprocal of G, except K7 which is a literal 0 ---
$ws.Range("K2").Formula = "=1/G2"
$ws.Range("K3:K6").Formula = "=1/G3"
$ws.Range("K7").Value = 0

# --- New header cell K1, same text as G1 ---
$ws.Range("K1").Value = $ws.Range("G1").Text

# --- Column width for new column K ---
$ws.Columns("K").ColumnWidth = 14

# --- Header alignment: center D1:E1 and H1:K1 ---
$ws.Range("D1:E1").HorizontalAlignment = -4108
$ws.Range("H1:K1").HorizontalAlignment = -4108

# --- Center alignment for D/E data rows 2-7 ---
$ws.Range("D2:E7").HorizontalAlignment = -4108

# --- Center alignment for H/I/J/K data rows 2-7 ---
$ws.Range("H2:K7").HorizontalAlignment = -4108

# --- Borders around D:E block (rows 1-7) ---
$ws.Range("D1:E7").BorderAround(1)
$ws.Range("D1:E7").Borders.Item(11).LineStyle = 1
$ws.Range("E1:E7").Borders.Item(10).LineStyle = 1

# --- Borders around H:K block (rows 1-7), right edge boxed ---
$ws.Range("H1:K7").BorderAround(1)
$ws.Range("H1:K7").Borders.Item(11).LineStyle = 1
$ws.Range("K1:K7").Borders.Item(10).LineStyle = 1

# --- Selection, as recorded after the edit ---
$ws.Range("J2:K7").Select()
